$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 250
$ws.Range("B2").Value = 93
$ws.Range("C2").Value = 938
$ws.Range("D2").Value = 24031
$ws.Range("E2").Value = 15625
$ws.Range("F2").Value = 219
$ws.Range("G2").Value = 219
$ws.Range("H2").Value = 313
$ws.Range("I2").Value = 4171
$ws.Range("J2").Value = 4016
$ws.Range("K2").Value = 48469
$ws.Range("L2").Value = 3559000
$ws.Range("M2").Value = 2110
$ws.Range("A3").Value = 266
$ws.Range("B3").Value = 79
$ws.Range("C3").Value = 937
$ws.Range("E3").Value = 15672
$ws.Range("F3").Value = 203
$ws.Range("G3").Value = 218
$ws.Range("H3").Value = 328
$ws.Range("I3").Value = 4204
$ws.Range("J3").Value = 4047
$ws.Range("K3").Value = 48437
$ws.Range("M3").Value = 2125
$ws.Range("A4").Value = 266
$ws.Range("B4").Value = 78
$ws.Range("C4").Value = 938
$ws.Range("D4").Value = 24078
$ws.Range("E4").Value = 15656
$ws.Range("F4").Value = 219
$ws.Range("G4").Value = 219
$ws.Range("H4").Value = 344
$ws.Range("I4").Value = 4171
$ws.Range("J4").Value = 4062
$ws.Range("K4").Value = 48438
$ws.Range("M4").Value = 2109
$ws.Range("A7").Value = 438
$ws.Range("B7").Value = 281
$ws.Range("C7").Value = 1578
$ws.Range("D7").Value = 61063
$ws.Range("E7").Value = 2053172
$ws.Range("F7").Value = 4687
$ws.Range("G7").Value = 219
$ws.Range("H7").Value = 313
$ws.Range("I7").Value = 5407
$ws.Range("J7").Value = 10437
$ws.Range("K7").Value = 222703
$ws.Range("L7").Value = 7793156
$ws.Range("M7").Value = 4693625
$ws.Range("A8").Value = 469
$ws.Range("B8").Value = 281
$ws.Range("C8").Value = 1531
$ws.Range("D8").Value = 60906
$ws.Range("F8").Value = 4688
$ws.Range("G8").Value = 219
$ws.Range("H8").Value = 313
$ws.Range("I8").Value = 5375
$ws.Range("J8").Value = 10453
$ws.Range("A9").Value = 453
$ws.Range("B9").Value = 281
$ws.Range("C9").Value = 1547
$ws.Range("D9").Value = 62594
$ws.Range("F9").Value = 4687
$ws.Range("G9").Value = 234
$ws.Range("H9").Value = 343
$ws.Range("I9").Value = 5375
$ws.Range("J9").Value = 10422

$ws.Range("L2").Select()
